# Add season-record columns (Wins / Losses / Ties) to the stats table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1): AD="Wins", AE="Losses", AF="Ties" ---
# Copy the formatting of the existing last header cell (AC1) onto the new
# header cells so they match the bold/bordered/centered header style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data (rows 2-44): each player's team season record ---
$lastRow = 44
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 71  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 91  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
